$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The SqlIP value in C2 was wrong for the Linux deployment (127.0.0.1 only
# resolves to localhost, not the actual DB host) - fix it to the real
# server address. Store it as text (matching the sheet's other "@" / text
# formatted id-like columns) so the IP string isn't coerced into a number.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "192.168.1.113"

# Column C needs to widen to fit the longer IP string.
$ws.Range("C1").ColumnWidth = 100 / 7

# The active selection after the fix is on the corrected cell.
$ws.Range("C2").Select()
